# Update the "Syndication DTCManualApi" sheet: re-sort the Vendor/Status
# rows alphabetically by vendor and drop the obsolete "Zomato" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Syndication DTCManualApi")

$ws.Range("B2").Value = "Apple"
$ws.Range("C2").Value = "In Progress"

$ws.Range("B3").Value = "Bing"
$ws.Range("C3").Value = "Submitted"

$ws.Range("B4").Value = "Facebook"
$ws.Range("C4").Value = "Submitted"

$ws.Range("B5").Value = "Factual"
$ws.Range("C5").Value = "Submitted"

$ws.Range("B6").Value = "Foursquare"
$ws.Range("C6").Value = "Submitted"

$ws.Range("B7").Value = "Google"
$ws.Range("C7").Value = "Submitted"

$ws.Range("B8").Value = "HERE"
$ws.Range("C8").Value = "In Progress"

$ws.Range("B9").Value = "TomTom"
$ws.Range("C9").Value = "In Progress"

$ws.Rows.Item(10).Delete()
